# Update column G ("K" - strikeouts) values for rows 2-22 on Sheet1,
# replacing the previously-stored "Strike#" derived values with the
# newly regenerated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 4
    3  = 5
    4  = 3
    5  = 4
    6  = 3
    7  = 7
    8  = 6
    9  = 2
    10 = 5
    11 = 3
    12 = 3
    13 = 3
    14 = 2
    15 = 2
    16 = 2
    17 = 6
    18 = 3
    19 = 6
    20 = 7
    21 = 7
    22 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
